$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.875.48"
$ws.Range("E2").Value = "  -0.70%  "
$ws.Range("D3").Value = "1.871.30"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'304.80"
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("D6").Value = "'1.0000"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "'0.5080"
$ws.Range("E7").Value = "  -1.59%  "
$ws.Range("D8").Value = "'0.3664"
$ws.Range("E8").Value = "  -2.76%  "
$ws.Range("E9").Value = "  +0.32%  "
$ws.Range("D10").Value = "'0.8926"
$ws.Range("E10").Value = "  +0.36%  "
$ws.Range("D11").Value = "'20.68"
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("E12").Value = "  -1.15%  "
$ws.Range("D13").Value = "1.882.57"
$ws.Range("E13").Value = "  +0.83%  "
$ws.Range("D14").Value = "'94.88"
$ws.Range("E14").Value = "  +5.78%  "
$ws.Range("D15").Value = "'5.228"
$ws.Range("E15").Value = "  -1.54%  "
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").Value = "'0.000008503"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("D19").Value = "'0.9998"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").Value = "26.930.23"
$ws.Range("E20").Value = "  -0.59%  "
$ws.Range("D21").Value = "'5.018"
$ws.Range("E21").Value = "  -0.26%  "
$ws.Range("D22").Value = "2.121.83"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "'10.37"
$ws.Range("E23").Value = "  -1.47%  "
$ws.Range("D24").Value = "'6.393"
$ws.Range("E24").Value = "  -1.14%  "
$ws.Range("D25").Value = "'148.39"
$ws.Range("E25").Value = "  +0.47%  "
$ws.Range("D26").Value = "'1.777"
$ws.Range("E26").Value = "  -3.32%  "
$ws.Range("D27").Value = "'17.88"
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("D28").Value = "'2.093"
$ws.Range("E28").Value = "  -0.19%  "
$ws.Range("D29").Value = "'113.42"
$ws.Range("E29").Value = "  +0.59%  "
$ws.Range("E30").Value = "  +0.72%  "
$ws.Range("D31").Value = "'4.728"
$ws.Range("E31").Value = "  +1.04%  "
$ws.Range("D32").Value = "'0.09149"
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("D33").Value = "'0.05066"
$ws.Range("E33").Value = "  -1.01%  "
$ws.Range("D34").Value = "'0.7485"
$ws.Range("E34").Value = "  +3.11%  "
$ws.Range("D35").Value = "'2.968"
$ws.Range("E35").Value = "  -3.26%  "
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").Value = "'3.232"
$ws.Range("E37").Value = "  +5.55%  "
$ws.Range("D38").Value = "'2.524"
$ws.Range("E38").Value = "  +1.04%  "
$ws.Range("D40").Value = "'0.5588"
$ws.Range("E40").Value = "  +4.84%  "
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("D42").Value = "'6.628"
$ws.Range("E42").Value = "  +1.94%  "
$ws.Range("D43").Value = "'115.99"
$ws.Range("E43").Value = "  -0.26%  "
$ws.Range("D44").Value = "'8.585"
$ws.Range("E44").Value = "  +3.47%  "
$ws.Range("D45").Value = "'0.1478"
$ws.Range("E45").Value = "  +0.84%  "
$ws.Range("D46").Value = "'0.4771"
$ws.Range("E46").Value = "  +2.94%  "
$ws.Range("D47").Value = "'0.9997"
$ws.Range("D48").Value = "'10.13"
$ws.Range("E48").Value = "  +1.37%  "
$ws.Range("D49").Value = "'1.566"
$ws.Range("E49").Value = "  -0.40%  "
$ws.Range("D50").Value = "'37.00"
$ws.Range("E50").Value = "  +1.20%  "
$ws.Range("D51").Value = "'63.17"
$ws.Range("E51").Value = "  -0.68%  "
